$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values for rows 2-10, columns D:J (Optimal/Fixed N Rate, Yield, GM, Profit Loss)

# Row 2 - Barley / Standard
$ws.Range("D2").Value = 131
$ws.Range("E2").Value = 131
$ws.Range("F2").Value = 4.88
$ws.Range("G2").Value = 4.88
$ws.Range("H2").Value = 1152.27
$ws.Range("I2").Value = 1152.27

# Row 3 - Barley / +20%
$ws.Range("D3").Value = 134.5
$ws.Range("E3").Value = 131
$ws.Range("F3").Value = 4.89
$ws.Range("G3").Value = 4.88
$ws.Range("H3").Value = 1478.66
$ws.Range("I3").Value = 1478.16
$ws.Range("J3").Value = 0.5

# Row 4 - Barley / -20%
$ws.Range("D4").Value = 123.9
$ws.Range("E4").Value = 131
$ws.Range("F4").Value = 4.85
$ws.Range("G4").Value = 4.88
$ws.Range("H4").Value = 827.14
$ws.Range("I4").Value = 826.39
$ws.Range("J4").Value = 0.75

# Row 5 - Wheat / Standard
$ws.Range("D5").Value = 88.5
$ws.Range("E5").Value = 88.5
$ws.Range("F5").Value = 3.98
$ws.Range("G5").Value = 3.98
$ws.Range("H5").Value = 1022.76
$ws.Range("I5").Value = 1022.76

# Row 6 - Wheat / +20%
$ws.Range("D6").Value = 92
$ws.Range("E6").Value = 88.5
$ws.Range("F6").Value = 3.99
$ws.Range("G6").Value = 3.98
$ws.Range("H6").Value = 1312.79
$ws.Range("I6").Value = 1312.55
$ws.Range("J6").Value = 0.24

# Row 7 - Wheat / -20%
$ws.Range("D7").Value = 81.40000000000001
$ws.Range("E7").Value = 88.5
$ws.Range("F7").Value = 3.95
$ws.Range("G7").Value = 3.98
$ws.Range("H7").Value = 733.83
$ws.Range("I7").Value = 732.97
$ws.Range("J7").Value = 0.86

# Row 8 - Canola / Standard
$ws.Range("D8").Value = 146.9
$ws.Range("E8").Value = 146.9
$ws.Range("F8").Value = 1.68
$ws.Range("G8").Value = 1.68
$ws.Range("H8").Value = 746.8200000000001
$ws.Range("I8").Value = 746.8200000000001

# Row 9 - Canola / +20%
$ws.Range("D9").Value = 157.5
$ws.Range("E9").Value = 146.9
$ws.Range("F9").Value = 1.7
$ws.Range("G9").Value = 1.68
$ws.Range("H9").Value = 1000.56
$ws.Range("I9").Value = 999.4400000000001
$ws.Range("J9").Value = 1.13

# Row 10 - Canola / -20%
$ws.Range("D10").Value = 131
$ws.Range("E10").Value = 146.9
$ws.Range("F10").Value = 1.66
$ws.Range("G10").Value = 1.68
$ws.Range("H10").Value = 496.1
$ws.Range("I10").Value = 494.2
$ws.Range("J10").Value = 1.89
